# Daily countdown update:
# Column D = total day count for a cycle, Column E = days remaining,
# Column F = cycle start date (stored as an integer yyyyMMdd).
# Each day: E decrements by 1. When a row's remaining count has reached
# its last day (E = 1), it rolls over into a fresh cycle: E resets back
# to the row's total (D) and the start date (F) advances by D days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - total days
    $eCell = $ws.Cells.Item($r, 5)   # column E - days remaining
    $fCell = $ws.Cells.Item($r, 6)   # column F - start date (yyyyMMdd)

    # NOTE: reading `.Value` through this host mis-resolves to an empty
    # variant (silently coerces to 0 in arithmetic) - use `.Value2` for
    # reliable reads. Writes work fine via either, `.Value2` used for
    # consistency.
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # Column F must be a clean 8-digit yyyyMMdd date. Some rows carry a
    # corrupted value (e.g. a stray 9-digit typo) - those are left
    # completely untouched, same as the source data.
    $fText = [string][int]$fVal
    if (-not ($fText -match '^\d{8}$')) {
        continue
    }

    $parsedDate = $null
    try {
        $parsedDate = [DateTime]::ParseExact($fText, "yyyyMMdd", $null)
    }
    catch {
        continue
    }

    if ($eVal -le 1) {
        # Last day reached -> roll into a new cycle.
        $eCell.Value2 = $dVal

        $newDate = $parsedDate.AddDays([double]$dVal)
        $fCell.Value2 = [int]$newDate.ToString("yyyyMMdd")
    }
    else {
        $eCell.Value2 = $eVal - 1
    }
}
